$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Revise Midterm One weights (row 17 = "Midterm 1")
$ws.Range("C17").Value = 30
$ws.Range("D17").Value = 10

$ws.Range("F34").Select()
